$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Single quote all the fields in redirect methods in all controllers*") {
        $p.Range.Delete()
        break
    }
}
